# ------------------------------------------------------------------
# Refresh the "cryptos" price/volume snapshot (GitHub Actions style
# scheduled update). Only the cells whose values actually moved are
# touched; everything else in the sheet is left exactly as-is.
#
# Column D ("Price") values that happen to look like plain numbers
# (e.g. 603.36) are written with a leading apostrophe so Excel keeps
# them as literal text, same as the workbook's original text-formatted
# price cells (values such as "69.394.43" already fail numeric
# parsing because of the extra thousands separator, so they do not
# need the apostrophe).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.394.43'
$ws.Range("E2").Value = '  +1.43%  '

# Row 3
$ws.Range("D3").Value = '3.883.79'
$ws.Range("E3").Value = '  +1.45%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '''603.36'
$ws.Range("E5").Value = '  +0.54%  '

# Row 6
$ws.Range("D6").Value = '''169.95'
$ws.Range("E6").Value = '  +4.21%  '

# Row 7
$ws.Range("D7").Value = '3.883.81'
$ws.Range("E7").Value = '  +1.54%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = '''0.534'
$ws.Range("E9").Value = '  +1.06%  '

# Row 10
$ws.Range("E10").Value = '  +1.44%  '

# Row 11
$ws.Range("E11").Value = '  +1.34%  '

# Row 12
$ws.Range("D12").Value = '''0.467'
$ws.Range("E12").Value = '  +1.89%  '

# Row 13
$ws.Range("E13").Value = '  +4.67%  '

# Row 14
$ws.Range("D14").Value = '''38.23'
$ws.Range("E14").Value = '  +4.04%  '

# Row 15
$ws.Range("D15").Value = '4.538.33'
$ws.Range("E15").Value = '  +1.55%  '

# Row 16
$ws.Range("D16").Value = '3.870.06'
$ws.Range("E16").Value = '  +1.67%  '

# Row 17
$ws.Range("D17").Value = '69.470.69'

# Row 18
$ws.Range("D18").Value = '''18.74'
$ws.Range("E18").Value = '  +9.77%  '

# Row 19
$ws.Range("E19").Value = '  +0.74%  '

# Row 20
$ws.Range("E20").Value = '  -0.57%  '

# Row 21
$ws.Range("D21").Value = '''11.02'
$ws.Range("E21").Value = '  -1.51%  '

# Row 22
$ws.Range("D22").Value = '''488.54'
$ws.Range("E22").Value = '  +0.62%  '

# Row 23
$ws.Range("D23").Value = '''0.741'
$ws.Range("E23").Value = '  +3.51%  '

# Row 24
$ws.Range("E24").Value = '  +3.98%  '

# Row 25
$ws.Range("D25").Value = '''85.29'
$ws.Range("E25").Value = '  +1.51%  '

# Row 26
$ws.Range("D26").Value = '''2.29'
$ws.Range("E26").Value = '  +2.56%  '

# Row 27
$ws.Range("D27").Value = '''12.36'
$ws.Range("E27").Value = '  +2.35%  '

# Row 28
$ws.Range("D28").Value = '''10.11'
$ws.Range("E28").Value = '  +1.36%  '

# Row 29
$ws.Range("E29").Value = '  +0.27%  '

# Row 30
$ws.Range("D30").Value = '''2.98'
$ws.Range("E30").Value = '  +1.05%  '

# Row 31
$ws.Range("D31").Value = '4.034.76'
$ws.Range("E31").Value = '  +1.36%  '

# Row 32
$ws.Range("E32").Value = '  +1.39%  '

# Row 33
$ws.Range("D33").Value = '''7.80'
$ws.Range("E33").Value = '  -0.31%  '

# Row 34
$ws.Range("D34").Value = '''31.88'
$ws.Range("E34").Value = '  +0.37%  '

# Row 35
$ws.Range("D35").Value = '3.852.26'
$ws.Range("E35").Value = '  +2.06%  '

# Row 36
$ws.Range("E36").Value = '  +0.00%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.143'
$ws.Range("E37").Value = '  +2.78%  '

# Row 38
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = '''3.40'
$ws.Range("E38").Value = '  +14.96%  '

# Row 39
$ws.Range("D39").Value = '''6.09'
$ws.Range("E39").Value = '  +3.94%  '

# Row 40
$ws.Range("E40").Value = '  +0.56%  '

# Row 41
$ws.Range("E41").Value = '  +0.02%  '

# Row 42
$ws.Range("E42").Value = '  +2.44%  '

# Row 43
$ws.Range("D43").Value = '''2.07'
$ws.Range("E43").Value = '  +4.63%  '

# Row 44
$ws.Range("D44").Value = '''435.25'
$ws.Range("E44").Value = '  +1.70%  '

# Row 45
$ws.Range("D45").Value = '''47.98'
$ws.Range("E45").Value = '  -1.01%  '

# Row 46
$ws.Range("D46").Value = '''8.68'
$ws.Range("E46").Value = '  +3.32%  '

# Row 48
$ws.Range("D48").Value = '''0.000275'
$ws.Range("E48").Value = '  +21.82%  '

# Row 49
$ws.Range("E49").Value = '  +2.51%  '

# Row 50
$ws.Range("D50").Value = '''40.18'
$ws.Range("E50").Value = '  +4.03%  '

# Row 51
$ws.Range("D51").Value = '''141.29'
$ws.Range("E51").Value = '  -0.96%  '
